$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "me dá uma cortesia ?cortesia "
$ws.Range("B40").Value = "Poxa quem sabe você ganha ,advinha quem são os 5 ganhadores de bingo de dezembro de 2024 "
$ws.Range("B40").WrapText = $true

$ws.Range("B43").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
